$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tournament")

# --- 1. Insert two new rows right after the header row for the new
#        competition-key / host-key loc entries (copied, not moved, from
#        the tourney-level data that code will read going forward). ---
$ws.Rows("2:3").Insert()

$ws.Range("A2").Value = "competition-key"
$ws.Range("B2").Value = "mens-copa-america"

$ws.Range("A3").Value = "host-key"
$ws.Range("B3").Value = "usa"

# --- 2. Append the 14 venue-key rows (copied from the #venues sheet)
#        after the existing venue.1..venue.14 rows (now rows 7..20). ---
$venueKeys = @(
  "us-atlanta-ga",
  "us-arlington-tx",
  "us-santa-clara-ca",
  "us-houston-tx",
  "us-miami-fl",
  "us-inglewood-ca",
  "us-kansas-city-ks",
  "us-east-rutherford-nj",
  "us-las-vegas-nv",
  "us-glendale-az",
  "us-orlando-fl",
  "us-austin-tx",
  "us-kansas-city-mo",
  "us-charlotte-nc"
)

$startRow = 21
for ($i = 0; $i -lt $venueKeys.Length; $i++) {
  $r = $startRow + $i
  $ws.Range("A$r").Value = "venue-key." + ($i + 1)
  $ws.Range("B$r").Value = $venueKeys[$i]
}

# --- 3. Grow the `tournament` table to cover the new rows. ---
$lo = $ws.ListObjects.Item("tournament")
$lo.Resize($ws.Range("A1:I34"))

# --- 4. Match the refreshed selection left behind by the edit (cursor
#        dropped back up to row 2 after the inserts). ---
$ws.Activate()
$ws.Range("A2:XFD3").Select()
